$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New issue row: bug report + status
$ws.Range("A3").Value = "Map API key not working for signed APK."
$ws.Range("B3").Value = "Inprogress"

# Date reported, stored as a real date serial formatted as a short date (mm-dd-yy -> numFmtId 14)
$ws.Range("D3").Value = 42305
$ws.Range("D3").NumberFormat = "mm-dd-yy"

# Best-fit the new Date column and move the selection the way the author left it
$ws.Columns.Item(4).AutoFit()
$ws.Range("C3").Select()
